# IST price update 2025-12-20 19:26
# Insert a new column before column B (shifting B:F -> C:G) and populate the
# new column B with a fresh timestamp header plus a duplicate snapshot of
# the prices that used to live in column B (now shifted into column C).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank column at B; everything from B.. shifts right to C..
$ws.Columns.Item(2).Insert()

# Excel's native column-insert doesn't carry over the custom width from the
# neighbour, so restore it explicitly (source col width="21" -> ColumnWidth ~20.17).
$ws.Columns.Item(2).ColumnWidth = 20.17

# Clone the header-row format (bordered/bold/centered style) from the
# neighbouring header cell onto the new one, then set its own value.
$ws.Cells.Item(1, 3).Copy()
$ws.Cells.Item(1, 2).PasteSpecial(-4122)
$ws.Cells.Item(1, 2).Value = "2025-12-21 00:52"

# The new column's prices mirror the prices that were in column B before the
# insert (same values, unchanged) for every data row.
$prices = @{
    2  = 929
    3  = 569
    4  = 299
    5  = 569
    6  = 499
    7  = 569
    8  = 929
    9  = 299
    10 = 299
    11 = 929
    12 = 569
    13 = 569
    14 = 499
    15 = 499
    16 = 299
    17 = 929
    18 = 499
    19 = 1497
    20 = 929
    21 = 499
    22 = 299
    23 = 1299
    24 = 929
    25 = 929
    26 = 1299
}

foreach ($r in $prices.Keys) {
    $ws.Cells.Item($r, 2).Value = $prices[$r]
}
